$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1, columns N-R (new columns)
$ws.Range("N1").Value = "effect_size"
$ws.Range("O1").Value = "effect_type"
$ws.Range("P1").Value = "ci_lower"
$ws.Range("Q1").Value = "ci_upper"
$ws.Range("R1").Value = "effect_size_for_plotting"

# Row 2
$ws.Range("N2").Value = 0.999987664451351
$ws.Range("O2").Value = "OR"
$ws.Range("P2").Value = 0.999136414139482
$ws.Range("Q2").Value = 1.00083964001663
$ws.Range("R2").Value = -0.0000123355486492116

# Row 3
$ws.Range("N3").Value = -0.0017861531232811
$ws.Range("O3").Value = "Coefficient"
$ws.Range("P3").Value = -0.00809521828287318
$ws.Range("Q3").Value = 0.00452291203631098
$ws.Range("R3").Value = -0.0017861531232811

# Row 4
$ws.Range("N4").Value = -0.0141370682902708
$ws.Range("O4").Value = "Coefficient"
$ws.Range("P4").Value = -0.0180936028629444
$ws.Range("Q4").Value = -0.0101805337175972
$ws.Range("R4").Value = -0.0141370682902708

# Row 5
$ws.Range("N5").Value = 0.0103747831351027
$ws.Range("O5").Value = "Coefficient"
$ws.Range("P5").Value = -0.00138398817569522
$ws.Range("Q5").Value = 0.0221335544459006
$ws.Range("R5").Value = 0.0103747831351027

# Row 6
$ws.Range("N6").Value = 0.0649197239967274
$ws.Range("O6").Value = "Coefficient"
$ws.Range("P6").Value = 0.00830681853288331
$ws.Range("Q6").Value = 0.121532629460571
$ws.Range("R6").Value = 0.0649197239967274

# Row 7
$ws.Range("N7").Value = 0.0363871019405857
$ws.Range("O7").Value = "Coefficient"
$ws.Range("P7").Value = 0.00324583916208586
$ws.Range("Q7").Value = 0.0695283647190855
$ws.Range("R7").Value = 0.0363871019405857

# Row 8
$ws.Range("N8").Value = 0.858334725689981
$ws.Range("O8").Value = "Coefficient"
$ws.Range("P8").Value = -0.022670264845207
$ws.Range("Q8").Value = 1.73933971622517
$ws.Range("R8").Value = 0.858334725689981

# Row 9
$ws.Range("N9").Value = 0.0874462453173206
$ws.Range("O9").Value = "Coefficient"
$ws.Range("P9").Value = -0.0435598033799183
$ws.Range("Q9").Value = 0.218452294014559
$ws.Range("R9").Value = 0.0874462453173206

# Row 10
$ws.Range("N10").Value = -0.014851412692442
$ws.Range("O10").Value = "Coefficient"
$ws.Range("P10").Value = -0.127033831276105
$ws.Range("Q10").Value = 0.0973310058912205
$ws.Range("R10").Value = -0.014851412692442

# Row 11
$ws.Range("N11").Value = 0.858334725689981
$ws.Range("O11").Value = "Coefficient"
$ws.Range("P11").Value = -0.022670264845207
$ws.Range("Q11").Value = 1.73933971622517
$ws.Range("R11").Value = 0.858334725689981

# Row 12
$ws.Range("N12").Value = 0.997729165136341
$ws.Range("O12").Value = "OR"
$ws.Range("P12").Value = 0.994583837028644
$ws.Range("Q12").Value = 1.00088444020732
$ws.Range("R12").Value = -0.00227083486365942

# Row 13
$ws.Range("N13").Value = 0.998360256120945
$ws.Range("O13").Value = "OR"
$ws.Range("P13").Value = 0.994500992842594
$ws.Range("Q13").Value = 1.00223449566695
$ws.Range("R13").Value = -0.00163974387905497

# Row 14
$ws.Range("N14").Value = 6423908.52220609
$ws.Range("O14").Value = "Coefficient"
$ws.Range("P14").Value = -3083444.69934574
$ws.Range("Q14").Value = 15931261.7437579
$ws.Range("R14").Value = 6423908.52220609

# Row 15
$ws.Range("N15").Value = 7181931.49896579
$ws.Range("O15").Value = "Coefficient"
$ws.Range("P15").Value = -5679925.0396275
$ws.Range("Q15").Value = 20043788.0375591
$ws.Range("R15").Value = 7181931.49896579

# Row 16
$ws.Range("N16").Value = 0.534795381854538
$ws.Range("O16").Value = "IRR"
$ws.Range("P16").Value = 0.230776364764813
$ws.Range("Q16").Value = 1.23932145626964
$ws.Range("R16").Value = -0.465204618145462

# Row 17
$ws.Range("N17").Value = 581860.58694433
$ws.Range("O17").Value = "IRR"
$ws.Range("P17").Value = 0.0169031092697272
$ws.Range("Q17").Value = 20029554162887.1
$ws.Range("R17").Value = 581859.58694433

# Row 18
$ws.Range("N18").Value = 0.999865392486336
$ws.Range("O18").Value = "IRR"
$ws.Range("P18").Value = 0.99967757162987
$ws.Range("Q18").Value = 1.00005324863085
$ws.Range("R18").Value = -0.000134607513664342

